$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply data updates (Price / Volume(1h) columns) from the latest cryptos refresh.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '59.320.49'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -1.59%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.578.61'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -1.87%  '

$ws.Range('E4').Value = '  -0.21%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '556.12'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.85%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '142.03'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -2.62%  '

$ws.Range('E7').Value = '  +0.13%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.599'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -1.42%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.586.08'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -2.46%  '

$ws.Range('E10').Value = '  -1.64%  '

$ws.Range('E11').Value = '  -0.80%  '

$ws.Range('E12').Value = '  +12.24%  '

$ws.Range('E13').Value = '  +2.05%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '3.034.28'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -2.79%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '59.311.81'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -1.70%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '22.97'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +4.40%  '

$ws.Range('E17').Value = '  -0.17%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.586.89'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -2.39%  '

$ws.Range('E19').Value = '  +0.28%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '337.61'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.33%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '10.34'
$ws.Range('D21').Style = 'Normal'

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.46'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +1.06%  '

$ws.Range('E23').Value = '  -0.07%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.478'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +9.07%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '62.47'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -5.23%  '

$ws.Range('E26').Value = '  +0.22%  '

$ws.Range('E27').Value = '  -3.05%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.40'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.33%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.0₃0775'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -2.92%  '

$ws.Range('E30').Value = '  -0.05%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.23'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +1.42%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.68'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -1.59%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '158.87'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.70%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '19.07'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.33%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.10'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.03%  '

$ws.Range('E36').Value = '  +1.47%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.896'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +1.28%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '37.38'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.62%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.852'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -4.26%  '

$ws.Range('E40').Value = '  -1.65%  '

$ws.Range('E41').Value = '  +1.13%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '289.47'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -2.91%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '137.15'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +7.58%  '

$ws.Range('E44').Value = '  +0.34%  '

$ws.Range('E45').Value = '  -1.04%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.593'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -1.69%  '

$ws.Range('E47').Value = '  -0.35%  '

$ws.Range('E48').Value = '  -2.42%  '

$ws.Range('E49').Value = '  -0.06%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '18.71'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.19%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.938.48'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -1.11%  '
